$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3610
$ws.Range("C3").Value = 3828
$ws.Range("C4").Value = 3921
$ws.Range("C5").Value = 4079
$ws.Range("C6").Value = 4118
$ws.Range("C7").Value = 4118
$ws.Range("C8").Value = 4118
$ws.Range("C9").Value = 4118
$ws.Range("C10").Value = 4118
$ws.Range("C11").Value = 4118
$ws.Range("C12").Value = 4511
